$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205262/1743535585915/aNjCWRFZTU-camera-video-segment-driver-1743535590915.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBIFHNCLP%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T210752Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJGMEQCIBrT7O9UVhgnj1q5MpD98ZALuNsqfe4R9q7ccbnsqnIVAiBO6ylfpHn6S4KG78WlAgxnCQg7%2FJGwSwtBKinrP5nz4SrmAwjb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM1VH3UKA7u%2BqCeQNOKroDB1sJ6Inqk1pIpz2i06%2BvxuFsC4b4gUerIpbEvMvhYssjBvQu6oqL8ikAq%2F8JlV3RTMKvsAHqT0F9s5sBGCVswuSU3k6kZHakrLHyxbF0ZOJ2oCUMPnSMhX8%2F9aovS7QCBRflm6rAHUYnDOJ%2FBjd7d%2F029TtfEzOVvEvNTvebc8yYZeFz2NxsvJC%2Fxctm6kcOs0J4vZ2QNQ5l09msYL8MwlWXsJSdX4SSMyrPGKMZqU0RPmpEbWv0SO6DPM57FUz%2FuwUiDBJSIvAyWbtVocj2drgVqkBzFDW50Bt9BP7oTpOv2d9JJSWTUi0Hzs6LiTzFxE3jNh%2BIs3Fjdg%2BFkithJQZXCE7axzvOYnkkeZ2nI%2F%2F1vR1KW98aIGY%2F8yLyGSPkZ%2B6Slov8htgFvgqk9E9P06X9RyTJ0YM8xoY7cQM9761l3RpBD7mbc3qN9kajbfnhzN1jb0UEBDNsYP3cLhXk9NZnGw4EQD7Osktbk8%2FJlkFe7PZ4MKdFiAV3ROOw4S1tnpv6vdrhIiWLC0quV6NRSll5zG9O3RwSpjIubUzmJPKNXu5E10W%2FlY9WugFxr65pY3md1BPj39wssTCu7rW%2FBjqmAYYC0x5Se41%2FRITA%2B8HyQMReOlCk1fEqzrjAPcNwDu%2FApsCz%2BPpy%2Br12DqFeQydQ14adnex1BO%2FsiV0Lj9F%2FONT%2BwNXLUmvaMovmnc8%2BNiZRTw0xb7l7ZszZeyoYW8TUI3jvYOZfqfoHlbJjga7N%2FQ47t0d%2FWRbNJ72P9dlYtcReOxw%2BmlaA%2FmkfvS29mmK4oxWTGu54Bqn0SeVdEyvbXo29xY7rrWA%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2005%3A07%3A52%20GMT&X-Amz-Signature=342123d283f592476f490302388508b7764b2f48594576581d4c94c39039b3b0'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535585915/3RxOLQEM5t-camera-video-segment-1743535590915.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBIFHNCLP%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T210752Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJGMEQCIBrT7O9UVhgnj1q5MpD98ZALuNsqfe4R9q7ccbnsqnIVAiBO6ylfpHn6S4KG78WlAgxnCQg7%2FJGwSwtBKinrP5nz4SrmAwjb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM1VH3UKA7u%2BqCeQNOKroDB1sJ6Inqk1pIpz2i06%2BvxuFsC4b4gUerIpbEvMvhYssjBvQu6oqL8ikAq%2F8JlV3RTMKvsAHqT0F9s5sBGCVswuSU3k6kZHakrLHyxbF0ZOJ2oCUMPnSMhX8%2F9aovS7QCBRflm6rAHUYnDOJ%2FBjd7d%2F029TtfEzOVvEvNTvebc8yYZeFz2NxsvJC%2Fxctm6kcOs0J4vZ2QNQ5l09msYL8MwlWXsJSdX4SSMyrPGKMZqU0RPmpEbWv0SO6DPM57FUz%2FuwUiDBJSIvAyWbtVocj2drgVqkBzFDW50Bt9BP7oTpOv2d9JJSWTUi0Hzs6LiTzFxE3jNh%2BIs3Fjdg%2BFkithJQZXCE7axzvOYnkkeZ2nI%2F%2F1vR1KW98aIGY%2F8yLyGSPkZ%2B6Slov8htgFvgqk9E9P06X9RyTJ0YM8xoY7cQM9761l3RpBD7mbc3qN9kajbfnhzN1jb0UEBDNsYP3cLhXk9NZnGw4EQD7Osktbk8%2FJlkFe7PZ4MKdFiAV3ROOw4S1tnpv6vdrhIiWLC0quV6NRSll5zG9O3RwSpjIubUzmJPKNXu5E10W%2FlY9WugFxr65pY3md1BPj39wssTCu7rW%2FBjqmAYYC0x5Se41%2FRITA%2B8HyQMReOlCk1fEqzrjAPcNwDu%2FApsCz%2BPpy%2Br12DqFeQydQ14adnex1BO%2FsiV0Lj9F%2FONT%2BwNXLUmvaMovmnc8%2BNiZRTw0xb7l7ZszZeyoYW8TUI3jvYOZfqfoHlbJjga7N%2FQ47t0d%2FWRbNJ72P9dlYtcReOxw%2BmlaA%2FmkfvS29mmK4oxWTGu54Bqn0SeVdEyvbXo29xY7rrWA%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2005%3A07%3A52%20GMT&X-Amz-Signature=0a25c20be64d9ff3fc8b08ea349bfa43b25b496edbc8b0c943e3e424fa395ac4'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535541388/Uf5dS3fS3z-camera-video-segment-driver-1743535543888.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBIFHNCLP%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T210752Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJGMEQCIBrT7O9UVhgnj1q5MpD98ZALuNsqfe4R9q7ccbnsqnIVAiBO6ylfpHn6S4KG78WlAgxnCQg7%2FJGwSwtBKinrP5nz4SrmAwjb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM1VH3UKA7u%2BqCeQNOKroDB1sJ6Inqk1pIpz2i06%2BvxuFsC4b4gUerIpbEvMvhYssjBvQu6oqL8ikAq%2F8JlV3RTMKvsAHqT0F9s5sBGCVswuSU3k6kZHakrLHyxbF0ZOJ2oCUMPnSMhX8%2F9aovS7QCBRflm6rAHUYnDOJ%2FBjd7d%2F029TtfEzOVvEvNTvebc8yYZeFz2NxsvJC%2Fxctm6kcOs0J4vZ2QNQ5l09msYL8MwlWXsJSdX4SSMyrPGKMZqU0RPmpEbWv0SO6DPM57FUz%2FuwUiDBJSIvAyWbtVocj2drgVqkBzFDW50Bt9BP7oTpOv2d9JJSWTUi0Hzs6LiTzFxE3jNh%2BIs3Fjdg%2BFkithJQZXCE7axzvOYnkkeZ2nI%2F%2F1vR1KW98aIGY%2F8yLyGSPkZ%2B6Slov8htgFvgqk9E9P06X9RyTJ0YM8xoY7cQM9761l3RpBD7mbc3qN9kajbfnhzN1jb0UEBDNsYP3cLhXk9NZnGw4EQD7Osktbk8%2FJlkFe7PZ4MKdFiAV3ROOw4S1tnpv6vdrhIiWLC0quV6NRSll5zG9O3RwSpjIubUzmJPKNXu5E10W%2FlY9WugFxr65pY3md1BPj39wssTCu7rW%2FBjqmAYYC0x5Se41%2FRITA%2B8HyQMReOlCk1fEqzrjAPcNwDu%2FApsCz%2BPpy%2Br12DqFeQydQ14adnex1BO%2FsiV0Lj9F%2FONT%2BwNXLUmvaMovmnc8%2BNiZRTw0xb7l7ZszZeyoYW8TUI3jvYOZfqfoHlbJjga7N%2FQ47t0d%2FWRbNJ72P9dlYtcReOxw%2BmlaA%2FmkfvS29mmK4oxWTGu54Bqn0SeVdEyvbXo29xY7rrWA%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2005%3A07%3A52%20GMT&X-Amz-Signature=2dc71e6ef8ea5bca1e550b0a41a188578f9e825d5cc368b04430a66fc78cb04f'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532925726/R4KzSHnOTc-camera-video-segment-driver-1743532928226.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBIFHNCLP%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T210752Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJGMEQCIBrT7O9UVhgnj1q5MpD98ZALuNsqfe4R9q7ccbnsqnIVAiBO6ylfpHn6S4KG78WlAgxnCQg7%2FJGwSwtBKinrP5nz4SrmAwjb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM1VH3UKA7u%2BqCeQNOKroDB1sJ6Inqk1pIpz2i06%2BvxuFsC4b4gUerIpbEvMvhYssjBvQu6oqL8ikAq%2F8JlV3RTMKvsAHqT0F9s5sBGCVswuSU3k6kZHakrLHyxbF0ZOJ2oCUMPnSMhX8%2F9aovS7QCBRflm6rAHUYnDOJ%2FBjd7d%2F029TtfEzOVvEvNTvebc8yYZeFz2NxsvJC%2Fxctm6kcOs0J4vZ2QNQ5l09msYL8MwlWXsJSdX4SSMyrPGKMZqU0RPmpEbWv0SO6DPM57FUz%2FuwUiDBJSIvAyWbtVocj2drgVqkBzFDW50Bt9BP7oTpOv2d9JJSWTUi0Hzs6LiTzFxE3jNh%2BIs3Fjdg%2BFkithJQZXCE7axzvOYnkkeZ2nI%2F%2F1vR1KW98aIGY%2F8yLyGSPkZ%2B6Slov8htgFvgqk9E9P06X9RyTJ0YM8xoY7cQM9761l3RpBD7mbc3qN9kajbfnhzN1jb0UEBDNsYP3cLhXk9NZnGw4EQD7Osktbk8%2FJlkFe7PZ4MKdFiAV3ROOw4S1tnpv6vdrhIiWLC0quV6NRSll5zG9O3RwSpjIubUzmJPKNXu5E10W%2FlY9WugFxr65pY3md1BPj39wssTCu7rW%2FBjqmAYYC0x5Se41%2FRITA%2B8HyQMReOlCk1fEqzrjAPcNwDu%2FApsCz%2BPpy%2Br12DqFeQydQ14adnex1BO%2FsiV0Lj9F%2FONT%2BwNXLUmvaMovmnc8%2BNiZRTw0xb7l7ZszZeyoYW8TUI3jvYOZfqfoHlbJjga7N%2FQ47t0d%2FWRbNJ72P9dlYtcReOxw%2BmlaA%2FmkfvS29mmK4oxWTGu54Bqn0SeVdEyvbXo29xY7rrWA%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2005%3A07%3A52%20GMT&X-Amz-Signature=bf9105290b727ffc7ced5ff21ef504bc4b88384170ba3b238a87cdcfd7d944ef'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532886131/8jJIZWs74D-camera-video-segment-driver-1743532888631.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSBIFHNCLP%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T210752Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJGMEQCIBrT7O9UVhgnj1q5MpD98ZALuNsqfe4R9q7ccbnsqnIVAiBO6ylfpHn6S4KG78WlAgxnCQg7%2FJGwSwtBKinrP5nz4SrmAwjb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM1VH3UKA7u%2BqCeQNOKroDB1sJ6Inqk1pIpz2i06%2BvxuFsC4b4gUerIpbEvMvhYssjBvQu6oqL8ikAq%2F8JlV3RTMKvsAHqT0F9s5sBGCVswuSU3k6kZHakrLHyxbF0ZOJ2oCUMPnSMhX8%2F9aovS7QCBRflm6rAHUYnDOJ%2FBjd7d%2F029TtfEzOVvEvNTvebc8yYZeFz2NxsvJC%2Fxctm6kcOs0J4vZ2QNQ5l09msYL8MwlWXsJSdX4SSMyrPGKMZqU0RPmpEbWv0SO6DPM57FUz%2FuwUiDBJSIvAyWbtVocj2drgVqkBzFDW50Bt9BP7oTpOv2d9JJSWTUi0Hzs6LiTzFxE3jNh%2BIs3Fjdg%2BFkithJQZXCE7axzvOYnkkeZ2nI%2F%2F1vR1KW98aIGY%2F8yLyGSPkZ%2B6Slov8htgFvgqk9E9P06X9RyTJ0YM8xoY7cQM9761l3RpBD7mbc3qN9kajbfnhzN1jb0UEBDNsYP3cLhXk9NZnGw4EQD7Osktbk8%2FJlkFe7PZ4MKdFiAV3ROOw4S1tnpv6vdrhIiWLC0quV6NRSll5zG9O3RwSpjIubUzmJPKNXu5E10W%2FlY9WugFxr65pY3md1BPj39wssTCu7rW%2FBjqmAYYC0x5Se41%2FRITA%2B8HyQMReOlCk1fEqzrjAPcNwDu%2FApsCz%2BPpy%2Br12DqFeQydQ14adnex1BO%2FsiV0Lj9F%2FONT%2BwNXLUmvaMovmnc8%2BNiZRTw0xb7l7ZszZeyoYW8TUI3jvYOZfqfoHlbJjga7N%2FQ47t0d%2FWRbNJ72P9dlYtcReOxw%2BmlaA%2FmkfvS29mmK4oxWTGu54Bqn0SeVdEyvbXo29xY7rrWA%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2005%3A07%3A52%20GMT&X-Amz-Signature=c001fde2d7de41379907edb6187b6a8765f50284530b4022d09068a2abeabdb1'
